$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "69.231.56"
$ws.Range("E2").Value = "  -0.41%  "

# Row 3
$ws.Range("D3").Value = "2.488.22"
$ws.Range("E3").Value = "  -1.04%  "

# Row 4
$ws.Range("E4").Value = "  -0.01%  "

# Row 5
$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = "567.54"
$cell.ClearFormats()
$ws.Range("E5").Value = "  -0.98%  "

# Row 6
$cell = $ws.Range("D6")
$cell.NumberFormat = "@"
$cell.Value = "165.61"
$cell.ClearFormats()
$ws.Range("E6").Value = "  -0.31%  "

# Row 7
$ws.Range("E7").Value = "  -0.04%  "

# Row 8
$ws.Range("E8").Value = "  -0.81%  "

# Row 9
$ws.Range("E9").Value = "  -0.67%  "

# Row 10
$ws.Range("E10").Value = "  -0.91%  "

# Row 11
$cell = $ws.Range("D11")
$cell.NumberFormat = "@"
$cell.Value = "0.347"
$cell.ClearFormats()
$ws.Range("E11").Value = "  -2.76%  "

# Row 12
$ws.Range("E12").Value = "  -0.74%  "

# Row 13
$ws.Range("D13").Value = "2.944.04"

# Row 14
$ws.Range("D14").Value = "69.139.82"
$ws.Range("E14").Value = "  -0.29%  "

# Row 15
$ws.Range("E15").Value = "  -1.02%  "

# Row 16
$ws.Range("E16").Value = "  -3.01%  "

# Row 17
$ws.Range("D17").Value = "2.509.71"
$ws.Range("E17").Value = "  +0.24%  "

# Row 18
$cell = $ws.Range("D18")
$cell.NumberFormat = "@"
$cell.Value = "11.15"
$cell.ClearFormats()
$ws.Range("E18").Value = "  -1.12%  "

# Row 19
$cell = $ws.Range("D19")
$cell.NumberFormat = "@"
$cell.Value = "352.34"
$cell.ClearFormats()
$ws.Range("E19").Value = "  +0.71%  "

# Row 20
$cell = $ws.Range("D20")
$cell.NumberFormat = "@"
$cell.Value = "7.34"
$cell.ClearFormats()
$ws.Range("E20").Value = "  -2.89%  "

# Row 21
$ws.Range("E21").Value = "  -0.41%  "

# Row 22
$ws.Range("E22").Value = "  -3.31%  "

# Row 23
$ws.Range("E23").Value = "  -0.06%  "

# Row 24
$cell = $ws.Range("D24")
$cell.NumberFormat = "@"
$cell.Value = "69.06"
$cell.ClearFormats()
$ws.Range("E24").Value = "  -1.73%  "

# Row 25
$cell = $ws.Range("D25")
$cell.NumberFormat = "@"
$cell.Value = "3.79"
$cell.ClearFormats()
$ws.Range("E25").Value = "  -3.46%  "

# Row 26
$ws.Range("D26").Value = "2.616.56"
$ws.Range("E26").Value = "  -1.07%  "

# Row 27
$cell = $ws.Range("D27")
$cell.NumberFormat = "@"
$cell.Value = "8.58"
$cell.ClearFormats()
$ws.Range("E27").Value = "  -3.76%  "

# Row 28
$cell = $ws.Range("D28")
$cell.NumberFormat = "@"
$cell.Value = "1.01"
$cell.ClearFormats()
$ws.Range("E28").Value = "  +0.47%  "

# Row 29
$ws.Range("D29").Value = "0.0₃0868"
$ws.Range("E29").Value = "  -2.51%  "

# Row 30
$cell = $ws.Range("D30")
$cell.NumberFormat = "@"
$cell.Value = "3.62"
$cell.ClearFormats()
$ws.Range("E30").Value = "  +139.28%  "

# Row 31
$cell = $ws.Range("D31")
$cell.NumberFormat = "@"
$cell.Value = "7.48"
$cell.ClearFormats()
$ws.Range("E31").Value = "  -4.14%  "

# Row 32
$ws.Range("B32").Value = "Bittensor"
$ws.Range("C32").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$cell = $ws.Range("D32")
$cell.NumberFormat = "@"
$cell.Value = "437.76"
$cell.ClearFormats()
$ws.Range("E32").Value = "  -5.45%  "

# Row 33
$ws.Range("B33").Value = "Fetch.AI"
$ws.Range("C33").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$cell = $ws.Range("D33")
$cell.NumberFormat = "@"
$cell.Value = "1.19"
$cell.ClearFormats()
$ws.Range("E33").Value = "  -3.95%  "

# Row 34
$ws.Range("E34").Value = "  +0.02%  "

# Row 35
$ws.Range("E35").Value = "  -1.35%  "

# Row 36
$cell = $ws.Range("D36")
$cell.NumberFormat = "@"
$cell.Value = "154.46"
$cell.ClearFormats()

# Row 37
$cell = $ws.Range("D37")
$cell.NumberFormat = "@"
$cell.Value = "0.113"
$cell.ClearFormats()
$ws.Range("E37").Value = "  -3.44%  "

# Row 38
$ws.Range("E38").Value = "  -0.27%  "

# Row 39
$cell = $ws.Range("D39")
$cell.NumberFormat = "@"
$cell.Value = "18.09"
$cell.ClearFormats()
$ws.Range("E39").Value = "  -2.10%  "

# Row 40
$ws.Range("E40").Value = "  +0.03%  "

# Row 41
$cell = $ws.Range("D41")
$cell.NumberFormat = "@"
$cell.Value = "0.312"
$cell.ClearFormats()
$ws.Range("E41").Value = "  -2.01%  "

# Row 42
$cell = $ws.Range("D42")
$cell.NumberFormat = "@"
$cell.Value = "4.57"
$cell.ClearFormats()
$ws.Range("E42").Value = "  -2.57%  "

# Row 43
$ws.Range("E43").Value = "  -2.12%  "

# Row 44
$ws.Range("E44").Value = "  -1.57%  "

# Row 45
$ws.Range("E45").Value = "  -4.65%  "

# Row 46
$cell = $ws.Range("D46")
$cell.NumberFormat = "@"
$cell.Value = "138.17"
$cell.ClearFormats()
$ws.Range("E46").Value = "  -2.53%  "

# Row 47
$ws.Range("E47").Value = "  -1.40%  "

# Row 48
$cell = $ws.Range("D48")
$cell.NumberFormat = "@"
$cell.Value = "0.503"
$cell.ClearFormats()
$ws.Range("E48").Value = "  -3.14%  "

# Row 49
$ws.Range("E49").Value = "  -1.43%  "

# Row 50
$ws.Range("E50").Value = "  -0.84%  "

# Row 51
$ws.Range("E51").Value = "  -0.47%  "
